$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 253, shifting the existing rows
# 253..271 down to 254..272 (preserving all their data/formatting).
$ws.Rows.Item(253).Insert()

# Populate the newly inserted row 253 with the new weekly record.
$ws.Cells.Item(253, 1).Value = 4
$ws.Cells.Item(253, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(253, 3).Value = "Los Lagos"
$ws.Cells.Item(253, 4).Value = 44585
$ws.Cells.Item(253, 5).Value = 10
$ws.Cells.Item(253, 6).Value = 100112023
$ws.Cells.Item(253, 7).Value = "Brócoli"
$ws.Cells.Item(253, 8).Value = "Sin especificar"
$ws.Cells.Item(253, 9).Value = "Primera"
$ws.Cells.Item(253, 10).Value = 250
$ws.Cells.Item(253, 11).Value = 1500
$ws.Cells.Item(253, 12).Value = 1500
$ws.Cells.Item(253, 13).Value = 1500
$ws.Cells.Item(253, 14).Value = "$/unidad"
$ws.Cells.Item(253, 15).Value = "Región Metropolitana"
$ws.Cells.Item(253, 16).Value = 1500
$ws.Cells.Item(253, 17).Value = 1
$ws.Cells.Item(253, 18).Value = "Hortaliza"
